$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("256:256").Insert()

$ws.Range("A256").Value = 4
$ws.Range("B256").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C256").Value = "Los Lagos"
$ws.Range("D256").Value = 44543
$ws.Range("E256").Value = 10
$ws.Range("F256").Value = 100112006
$ws.Range("G256").Value = "Repollo"
$ws.Range("H256").Value = "Crespo record"
$ws.Range("I256").Value = "Primera"
$ws.Range("J256").Value = 500
$ws.Range("K256").Value = 1200
$ws.Range("L256").Value = 1200
$ws.Range("M256").Value = 1200
$ws.Range("N256").Value = "`$/unidad"
$ws.Range("O256").Value = "Región del Maule"
$ws.Range("P256").Value = 1200
$ws.Range("Q256").Value = 1
$ws.Range("R256").Value = "Hortaliza"
